# Add 3 new quest rows (15-17) for the "Wondering Merchant" NPC quest chain,
# and widen columns A (name) and B (npc_id) to fit the new, longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: That's One Creepy Doll
$ws.Range("A15").Value = "That's One Creepy Doll"
$ws.Range("B15").Value = "Wondering Merchant"
$ws.Range("C15").Value = "Creepy Baby Doll"
$ws.Range("D15").Value = 30000
$ws.Range("G15").Value = "Shadow Plane Grimoire"

# Row 16: Beauty is in the eye of the beholder
$ws.Range("A16").Value = "Beauty is in the eye of the beholder"
$ws.Range("B16").Value = "Wondering Merchant"
$ws.Range("C16").Value = "Shadow Plane Grimoire"
$ws.Range("D16").Value = 50000
$ws.Range("G16").Value = "Eye of the Beholder"

# Row 17: The Return of the King's Crown
$ws.Range("A17").Value = "The Return of the King's Crown"
$ws.Range("B17").Value = "Wondering Merchant"
$ws.Range("C17").Value = "Eye of the Beholder"
$ws.Range("D17").Value = 75000
$ws.Range("G17").Value = "Dead Kings Crown"

# Widen column A (29 -> 43) and column B (17 -> 22) to fit the new text.
# (ColumnWidth round-trips through the engine with a fixed +5/6 padding, so
# subtract it here to land exactly on the target stored width.)
$ws.Columns.Item(1).ColumnWidth = 43 - 5/6
$ws.Columns.Item(2).ColumnWidth = 22 - 5/6
